# Insert a new weekly Puerro price record at row 273 of the single sheet,
# pushing the existing rows 273-286 down to 274-287 (dimension grows to R287).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 273 (shifts rows 273:286 -> 274:287).
$ws.Rows.Item(273).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A273").Value = 10
$ws.Range("B273").Value = "Vega Modelo de Temuco"
$ws.Range("C273").Value = "La Araucanía"
$ws.Range("D273").Value = 45041
$ws.Range("E273").Value = 9
$ws.Range("F273").Value = 100112005
$ws.Range("G273").Value = "Puerro"
$ws.Range("H273").Value = "Azul de Maquehue"
$ws.Range("I273").Value = "Primera"
$ws.Range("J273").Value = 40
$ws.Range("K273").Value = 11000
$ws.Range("L273").Value = 11000
$ws.Range("M273").Value = 11000
$ws.Range("N273").Value = "$/docena de paquetes"
$ws.Range("O273").Value = "Provincia de Cautín"
$ws.Range("P273").Value = 917
$ws.Range("Q273").Value = 12
$ws.Range("R273").Value = "Hortaliza"
